$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Address column entries for the two data rows.
# This causes the two now-unused shared strings ("Ripley DE5" and the
# full "Collins Earthworks Limited..." address) to be dropped from
# sharedStrings.xml, and all later shared-string indices shift down,
# matching the target workbook.
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()

# Move/update the active selection to F3 (was D2).
$ws.Range("F3").Select()
